$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.531202554702759
$ws.Range("B1").Value = 3.119097709655762
$ws.Range("C1").Value = 4.684049606323242
$ws.Range("D1").Value = 1.855604887008667
$ws.Range("E1").Value = 1.171880602836609
